# The deck's in-use theme (ppt/theme/theme2.xml, wired to the slide master /
# presentation) currently carries the "Integral" / "Red Violet" colour
# scheme, while the otherwise-unused ppt/theme/theme1.xml (wired only to the
# notes master) carries the default "Office Theme" / "Office" colour scheme.
# The authored edit swaps those two themes' contents so the presentation now
# renders with the stock Office colour scheme.
#
# The PowerPoint object model exposes the live theme's 12-colour scheme via
# Slide.ThemeColorScheme (ThemeColorScheme.Colors(i).RGB, COM "BBGGRR" long
# values, in dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink order) so we drive the
# swap through that rather than touching package parts directly.

function HexToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Target palette = the "Office" colour scheme (formerly theme1.xml), applied
# in clrScheme order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeHex = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = HexToComRgb $officeHex[$i - 1]
}
